$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new date, hours, and activity entry
$ws.Range("A17").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 42822
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Revising - Looking over project, learning some new technologies"

# Row 19: new date, hours, and activity entry
$ws.Range("A17").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 42827
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Learning - Learning how to use the tools that were selected"

# Update the active selection to A20
$ws.Range("A20").Select()
